# USD transactions suite added
#
# Updates row 1 of the "Transactions" sheet from an AUD/MasterCard/Tester
# record to a USD/Visa/InterTester record, and appends a new row 2 with a
# second USD/MasterCard/InterTester record.
#
# All six columns (A:F) are plain text-typed cells (shared strings) holding
# numeric-looking values (amounts, codes, big ids) as well as real text
# (tester name, currency, card brand) - so every write below forces a Text
# number format before assigning the value (otherwise Excel's automatic
# type-detection would silently convert numeric-looking strings like
# "125" or "3320605440" into numeric cells) and then restores the cell's
# style to "Normal" afterwards so the persisted cell keeps the workbook's
# default (unstyled) look, matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$columns = @("A", "B", "C", "D", "E", "F")

# Row 1: existing AUD/MasterCard/Tester record -> USD/Visa/InterTester record
$row1Values = @("InterTester", "USD", "2109.02", "125", "Visa", "3320605440")
for ($i = 0; $i -lt $columns.Length; $i++) {
    $cell = $ws.Range($columns[$i] + "1")
    Set-TextValue $cell $row1Values[$i]
}

# Row 2: new USD/MasterCard/InterTester record
$row2Values = @("InterTester", "USD", "2240.27", "125", "MasterCard", "3322910208")
for ($i = 0; $i -lt $columns.Length; $i++) {
    $cell = $ws.Range($columns[$i] + "2")
    Set-TextValue $cell $row2Values[$i]
}
